$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Step 1: merge the "16" / "." / " What is the WHERE clause ..." runs into
# a single run (same visible text, same bold/sz28 formatting of the first
# run) by doing a Find/Replace across the whole phrase.
# -----------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute(
    "16. What is the WHERE clause used for, and how is it used to filter data?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "16. What is the WHERE clause used for, and how is it used to filter data?",
    2) | Out-Null

# -----------------------------------------------------------------------
# Step 2: locate the "In SQL, the WHERE clause ..." paragraph (it stays
# unchanged) and insert the new Q17/A17 paragraphs right after it.
# -----------------------------------------------------------------------
$wherePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "In SQL, the WHERE clause is used to filter rows") {
        $wherePara = $p
    }
}

$wherePara.Range.InsertParagraphAfter() | Out-Null
$q17 = $wherePara.Next()

$q17xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
          '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>17</w:t></w:r>' +
          '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
          '<w:t xml:space="preserve"> How do you retrieve distinct values from a column in SQL? </w:t></w:r>' +
          '</w:p>'
$q17.Range.InsertXML($q17xml) | Out-Null

# Re-find the Q17 paragraph (ranges have shifted after InsertXML) so we can
# append the answer paragraph right after it.
$q17b = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "How do you retrieve distinct values from a column in SQL") {
        $q17b = $p
    }
}

$q17b.Range.InsertParagraphAfter() | Out-Null
$a17 = $q17b.Next()
$a17xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:r><w:t>Using the DISTINCT keyword in combination with the SELECT command, ' +
          'we can extract distinct values from a column in SQL. By filtering out duplicate ' +
          'values and returning only unique values from the specified column, the DISTINCT ' +
          'keyword is used</w:t></w:r></w:p>'
$a17.Range.InsertXML($a17xml) | Out-Null

Write-Output "done"
